$wb = $excel.ActiveWorkbook

# The workbook has two worksheets ("展览" and "全部类型") that share the
# same "想去人数" (column F) data for a number of rows. Update both in sync.
$sheetNames = @("展览", "全部类型")

$updates = @{
    2  = 1625
    3  = 9016
    7  = 432
    8  = 185
    11 = 3848
    13 = 380
    15 = 4292
    24 = 2654
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
